$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(288, 3568.78, 0.028020799999999998, 39.14, 4.3499999999999996, 32.07, 22.52, 1.02, 0.9),
    @(336, 3781.1210000000001, 0.026447200000000001, 35.51, 4.05, 36.11, 21.82, 1.54, 0.98),
    @(384, 4025.6379999999999, 0.0248408, 32.85, 3.89, 37.869999999999997, 23.52, 0.95, 0.92)
)

$startRow = 16
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}

$ws.Range("J18").Select()
